$d = $word.ActiveDocument

# Extracts the <w:rPr>...</w:rPr> (if any) that is directly set on the run
# whose <w:t> text equals $needleText, by scanning the WordOpenXML of the
# paragraph/range that contains it. This avoids relying on Range.Font.Bold /
# Range.Font.Italic, which reflect *inherited* (style-based) formatting
# rather than only the run's own explicit formatting.
function Get-RunRPr($wordOpenXml, $needleText) {
    $pattern = '<w:r>(.*?)<w:t[^>]*>([^<]*)</w:t>(.*?)</w:r>'
    $regex = [regex]$pattern
    $m = $regex.Match($wordOpenXml)
    while ($m.Success) {
        if ($m.Groups[2].Value -eq $needleText) {
            $inner = $m.Groups[1].Value
            if ($inner -match '(<w:rPr>.*?</w:rPr>)') {
                return $Matches[1]
            }
            return ""
        }
        $m = $m.NextMatch()
    }
    return ""
}

# Replaces exactly one occurrence of $searchText (an entire run's text) with
# $replaceText, preserving the run's own explicit formatting (rPr) as well as
# any sibling runs (e.g. the leading empty <w:r/> runs used in this document)
# by using InsertXML scoped to just the matched run instead of rewriting the
# whole paragraph the way Find/Replace or Range.Text assignment would.
function Replace-ExactRunOnce($searchText, $replaceText) {
    $probe = $d.Content
    $found = $probe.Find.Execute($searchText, $true, $true, $false, $false, $false, $true)
    if (-not $found) {
        return $false
    }

    $start = $probe.Start
    $end = $probe.End
    $rPr = Get-RunRPr $probe.WordOpenXML $searchText

    $escaped = $replaceText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $needsPreserve = ($replaceText -ne $replaceText.Trim())
    $spaceAttr = ""
    if ($needsPreserve) { $spaceAttr = ' xml:space="preserve"' }

    $runXml = "<w:r>" + $rPr + "<w:t$spaceAttr>" + $escaped + "</w:t></w:r>"
    $packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Re-fetch a fresh Range object at the saved offsets: reusing the Range
    # mutated by Find.Execute for InsertXML can insert alongside stale content
    # instead of replacing it.
    $target = $d.Range($start, $end)
    $target.InsertXML($packageXml)
    return $true
}

# Repeats the single replacement until the search text can no longer be
# found, so every occurrence in the document gets updated.
function Replace-AllRuns($searchText, $replaceText) {
    $guard = 0
    while ((Replace-ExactRunOnce $searchText $replaceText) -and ($guard -lt 25)) {
        $guard = $guard + 1
    }
}

# Title (Heading1) and the bold title run near the end share the same old
# text and both get replaced with the same new text.
Replace-AllRuns "Play La Messicana Slot for Free - Review 2021" "Play La Messicana for Free - Exciting Mexican Slot Game"

# "What we like" bullet points
Replace-AllRuns "Mexican-themed with typical elements and colors" "Authentic Mexican theme"
Replace-AllRuns "Good illustrations that add to the gaming experience" "Beautiful illustrations"
Replace-AllRuns "Features a Wild and a Scatter symbol" "Exciting gameplay mechanics"
Replace-AllRuns "Medium volatile with an RTP of 95.08%" "Chance to win random prizes"

# "What we don't like" bullet points
Replace-AllRuns "Low-value winnings could lead to an unnoticed drop in the game's balance" "Low-value winnings"
Replace-AllRuns "Melody could be overshadowed by pronounced sound effects" "Pronounced sound effects"

# Meta description (italic text)
Replace-AllRuns "Read our honest review of La Messicana by Cristaltec. Play this Mexican-themed slot for free, including pros and cons and similar slots suggestions." "Read our review of La Messicana and play this exciting Mexican-themed slot game for free."
